# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# Re-writes the worker data table (rows 16-22) on Hoja1 with the updated
# employee list / order coming from the refreshed database export, and
# updates the "Valor Mora" / "Salario Basico" amounts that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data for rows 16-22: Tipo Doc, N Doc, Nombre, Periodo Mora, Valor Mora, Salario Basico
$data = @(
    @(16, "CC", "73568861",   "ROBINSON HENAO CUADRADO",               "1804", 29509, 737717),
    @(17, "CC", "1043971127", "MARIA ISABEL GARCIA BERMUDEZ",          "1804", 60000, 1500000),
    @(18, "CC", "1143375774", "WILFRAM CARDOZA CARDOZA",               "1804", 29509, 737717),
    @(19, "CC", "73192420",   "HENRY JUSTIN CASTILLO MARTINEZ",        "1804", 29509, 737717),
    @(20, "CC", "1005675244", "YORCY JOSE ARROYO OSORIO",              "1804", 29509, 737717),
    @(21, "CC", "71295087",   "CRISTHIAN ALEXANDER CARTAGENA GIRALDO", "1804", 29509, 737717),
    @(22, "CC", "1052085622", "MARIA FERNANDA ATEHORTUA JARAMILLO",    "1804", 29509, 737717)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
